$d = $word.ActiveDocument

$d.Content.Find.Execute("**ID__AFFARS_mp_5315_4_topic_3__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_MP5315_4_2__ID**", 2)

$p = $d.Paragraphs.Item(1)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25
